# Updated setting_analysis script to include additional analyses
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column A (switch used) relabels for rows 5-16 ---
$ws.Range("A5").Value = "MX2A Brown"
$ws.Range("A6").Value = "MX2A Black"
$ws.Range("A7").Value = "MX Red"
$ws.Range("A8").Value = "MX Red"
$ws.Range("A9").Value = "MX2A Black"
$ws.Range("A10").Value = "MX2A Brown"
$ws.Range("A11").Value = "MX2A Brown"
$ws.Range("A12").Value = "MX2A Black"
$ws.Range("A13").Value = "MX Red"
$ws.Range("A14").Value = "MX Red"
$ws.Range("A15").Value = "MX2A Black"
$ws.Range("A16").Value = "MX2A Brown"

# --- Column D (average WPM values) updates ---
$ws.Range("D5").Value = 140.06100000000001
$ws.Range("D6").Value = 144.89599999999999
$ws.Range("D7").Value = 146.749
$ws.Range("D8").Value = 145.95400000000001
$ws.Range("D9").Value = 164.17
$ws.Range("D10").Value = 150.68100000000001
$ws.Range("D11").Value = 140.142
$ws.Range("D12").Value = 145.62899999999999
$ws.Range("D13").Value = 144.98099999999999
$ws.Range("D14").Value = 154.321
$ws.Range("D15").Value = 136.476
$ws.Range("D16").Value = 132.02000000000001

# --- Remove the old "cycle order" / "use for next set of tests" block (F9, J9, J11) ---
$ws.Range("F9").Style = "Normal"
$ws.Range("F9").Value = 1
$ws.Range("J9").ClearContents()
$ws.Range("J11").ClearContents()
$ws.Rows.Item(9).AutoFit()

# row 7 now carries the "Switches by cycle order:" heading
$ws.Range("F7").Value = "Switches by cycle order:"

# row 8 carries the "Y" flag + its explanatory note
$ws.Range("H8").Value = "Y"
$ws.Range("I8").Value = "A 'Y' in the columns to the left shows which cycle to do next."

# row 9 becomes the 1/2/3 column headers (plain numbers, not text)
$ws.Range("G9").Value = 2
$ws.Range("H9").Value = 3

# --- Cycle-order grid, rows 10-21, columns F:H ---
$ws.Range("F10").Value = "MX Red"
$ws.Range("G10").Value = "MX2A Brown"
$ws.Range("H10").Value = "MX2A Black"

$ws.Range("F11").Value = "MX2A Brown"
$ws.Range("G11").Value = "MX2A Black"
$ws.Range("H11").Value = "MX Red"

$ws.Range("F12").Value = "MX2A Black"
$ws.Range("G12").Value = "MX Red"
$ws.Range("H12").Value = "MX2A Brown"

$ws.Range("F13").Value = "MX2A Black"
$ws.Range("G13").Value = "MX Red"
$ws.Range("H13").Value = "MX2A Brown"

$ws.Range("F14").Value = "MX2A Brown"
$ws.Range("G14").Value = "MX2A Black"
$ws.Range("H14").Value = "MX Red"

$ws.Range("F15").Value = "MX Red"
$ws.Range("G15").Value = "MX2A Brown"
$ws.Range("H15").Value = "MX2A Black"

$ws.Range("F16").Value = "MX Red"
$ws.Range("G16").Value = "MX2A Brown"
$ws.Range("H16").Value = "MX2A Black"

$ws.Range("F17").Value = "MX2A Brown"
$ws.Range("G17").Value = "MX2A Black"
$ws.Range("H17").Value = "MX Red"

$ws.Range("F18").Value = "MX2A Black"
$ws.Range("G18").Value = "MX Red"
$ws.Range("H18").Value = "MX2A Brown"

$ws.Range("F19").Value = "MX2A Black"
$ws.Range("G19").Value = "MX Red"
$ws.Range("H19").Value = "MX2A Brown"

$ws.Range("F20").Value = "MX2A Brown"
$ws.Range("G20").Value = "MX2A Black"
$ws.Range("H20").Value = "MX Red"

$ws.Range("F21").Value = "MX Red"
$ws.Range("G21").Value = "MX2A Brown"
$ws.Range("H21").Value = "MX2A Black"

# --- "Use for next set of tests:" + Y flag, now moved to row 22 ---
$ws.Range("E22").Value = "Use for next set of tests:"
$ws.Range("G22").Value = "Y"

# --- Explanatory note, now on row 24 ---
$ws.Range("F24").Value = "(Note: if running additional sets of tests right after the first set of tests, use the same cycle that you did previously. For instance, if running two sets of tests back to back, cycle through your switches as follows: Red Brown Black Black Brown Red Red Brown Black Black Brown Red Red Brown Black Black Brown Red Red Brown Black Black Brown Red.)"

# --- Column width tweaks (G/H get explicit widths like column A) ---
$ws.Columns.Item(7).ColumnWidth = 12.3
$ws.Columns.Item(8).ColumnWidth = 11.17

# --- View state: scroll position + active selection ---
$win = $excel.ActiveWindow
$win.ScrollRow = 4
$win.ScrollColumn = 1
$ws.Range("E8").Select()
